# Fix mixed-up fixture rows in the "Serbia Prva Liga" results table.
# Several adjacent match rows had their data (id, teams, score, odds, ...)
# swapped between rows; this restores the correct per-row values while
# leaving the row rank (col A) and Div/Div Original Name/Date (cols C:E) intact.
$ws = $excel.ActiveWorkbook.ActiveSheet

# Row 53
$ws.Cells.Item(53, 2).Value2 = 6989506  # B53
$ws.Cells.Item(53, 6).Value2 = 'FK Radnicki Beograd'  # F53
$ws.Cells.Item(53, 7).Value2 = 'Radnicki Sremska Mitrovica'  # G53
$ws.Cells.Item(53, 8).Value2 = 0  # H53
$ws.Cells.Item(53, 9).Value2 = 1  # I53
$ws.Cells.Item(53, 10).Value2 = 'A'  # J53
$ws.Cells.Item(53, 11).Value2 = 2.5  # K53
$ws.Cells.Item(53, 13).Value2 = 2.625  # M53
$ws.Cells.Item(53, 14).Value2 = 2.5  # N53
$ws.Cells.Item(53, 15).Value2 = 3  # O53
$ws.Cells.Item(53, 16).Value2 = 2.625  # P53
$ws.Cells.Item(53, 17).Value2 = 0  # Q53
$ws.Cells.Item(53, 18).Value2 = 1.85  # R53
$ws.Cells.Item(53, 19).Value2 = 1.95  # S53
$ws.Cells.Item(53, 21).Value2 = 2  # U53
$ws.Cells.Item(53, 22).Value2 = 1.8  # V53
$ws.Cells.Item(53, 23).Value2 = -1  # W53
$ws.Cells.Item(53, 25).Value2 = 1.625  # Y53
$ws.Cells.Item(53, 26).Value2 = -1  # Z53
$ws.Cells.Item(53, 27).Value2 = 0.95  # AA53
$ws.Cells.Item(53, 29).Value2 = 0.8  # AC53

# Row 54
$ws.Cells.Item(54, 2).Value2 = 6989507  # B54
$ws.Cells.Item(54, 6).Value2 = 'RFK Novi Sad 1921'  # F54
$ws.Cells.Item(54, 7).Value2 = 'FK Tekstilac Odzaci'  # G54
$ws.Cells.Item(54, 8).Value2 = 1  # H54
$ws.Cells.Item(54, 9).Value2 = 0  # I54
$ws.Cells.Item(54, 10).Value2 = 'H'  # J54
$ws.Cells.Item(54, 11).Value2 = 3  # K54
$ws.Cells.Item(54, 13).Value2 = 2.25  # M54
$ws.Cells.Item(54, 14).Value2 = 4  # N54
$ws.Cells.Item(54, 15).Value2 = 3.2  # O54
$ws.Cells.Item(54, 16).Value2 = 1.8  # P54
$ws.Cells.Item(54, 17).Value2 = 0.5  # Q54
$ws.Cells.Item(54, 18).Value2 = 1.9  # R54
$ws.Cells.Item(54, 19).Value2 = 1.9  # S54
$ws.Cells.Item(54, 21).Value2 = 1.9  # U54
$ws.Cells.Item(54, 22).Value2 = 1.9  # V54
$ws.Cells.Item(54, 23).Value2 = 3  # W54
$ws.Cells.Item(54, 25).Value2 = -1  # Y54
$ws.Cells.Item(54, 26).Value2 = 0.8999999999999999  # Z54
$ws.Cells.Item(54, 27).Value2 = -1  # AA54
$ws.Cells.Item(54, 29).Value2 = 0.8999999999999999  # AC54

# Row 87
$ws.Cells.Item(87, 2).Value2 = 6989515  # B87
$ws.Cells.Item(87, 6).Value2 = 'OFK Vrsac'  # F87
$ws.Cells.Item(87, 7).Value2 = 'RFK Novi Sad 1921'  # G87
$ws.Cells.Item(87, 11).Value2 = 1.5  # K87
$ws.Cells.Item(87, 12).Value2 = 3.75  # L87
$ws.Cells.Item(87, 13).Value2 = 5.5  # M87
$ws.Cells.Item(87, 14).Value2 = 1.5  # N87
$ws.Cells.Item(87, 15).Value2 = 3.75  # O87
$ws.Cells.Item(87, 16).Value2 = 6  # P87
$ws.Cells.Item(87, 17).Value2 = -1  # Q87
$ws.Cells.Item(87, 18).Value2 = 1.825  # R87
$ws.Cells.Item(87, 19).Value2 = 1.975  # S87
$ws.Cells.Item(87, 20).Value2 = 2.25  # T87
$ws.Cells.Item(87, 21).Value2 = 2  # U87
$ws.Cells.Item(87, 22).Value2 = 1.8  # V87
$ws.Cells.Item(87, 23).Value2 = 0.5  # W87
$ws.Cells.Item(87, 26).Value2 = 0  # Z87
$ws.Cells.Item(87, 27).Value2 = -0  # AA87
$ws.Cells.Item(87, 29).Value2 = 0.8  # AC87

# Row 89
$ws.Cells.Item(89, 2).Value2 = 6989684  # B89
$ws.Cells.Item(89, 6).Value2 = 'Radnicki Sremska Mitrovica'  # F89
$ws.Cells.Item(89, 7).Value2 = 'FK Tekstilac Odzaci'  # G89
$ws.Cells.Item(89, 11).Value2 = 2.4  # K89
$ws.Cells.Item(89, 12).Value2 = 2.8  # L89
$ws.Cells.Item(89, 13).Value2 = 2.9  # M89
$ws.Cells.Item(89, 14).Value2 = 3.75  # N89
$ws.Cells.Item(89, 15).Value2 = 2.9  # O89
$ws.Cells.Item(89, 16).Value2 = 1.95  # P89
$ws.Cells.Item(89, 17).Value2 = 0.5  # Q89
$ws.Cells.Item(89, 18).Value2 = 1.8  # R89
$ws.Cells.Item(89, 19).Value2 = 2  # S89
$ws.Cells.Item(89, 20).Value2 = 2  # T89
$ws.Cells.Item(89, 21).Value2 = 1.85  # U89
$ws.Cells.Item(89, 22).Value2 = 1.95  # V89
$ws.Cells.Item(89, 23).Value2 = 2.75  # W89
$ws.Cells.Item(89, 26).Value2 = 0.8  # Z89
$ws.Cells.Item(89, 27).Value2 = -1  # AA89
$ws.Cells.Item(89, 29).Value2 = 0.95  # AC89

# Row 112
$ws.Cells.Item(112, 2).Value2 = 6989691  # B112
$ws.Cells.Item(112, 6).Value2 = 'Radnicki Sremska Mitrovica'  # F112
$ws.Cells.Item(112, 7).Value2 = 'Smederevo'  # G112
$ws.Cells.Item(112, 8).Value2 = 0  # H112
$ws.Cells.Item(112, 9).Value2 = 0  # I112
$ws.Cells.Item(112, 11).Value2 = 2  # K112
$ws.Cells.Item(112, 13).Value2 = 3.6  # M112
$ws.Cells.Item(112, 14).Value2 = 1.95  # N112
$ws.Cells.Item(112, 15).Value2 = 3  # O112
$ws.Cells.Item(112, 16).Value2 = 3.75  # P112
$ws.Cells.Item(112, 17).Value2 = -0.5  # Q112
$ws.Cells.Item(112, 18).Value2 = 2  # R112
$ws.Cells.Item(112, 19).Value2 = 1.8  # S112
$ws.Cells.Item(112, 21).Value2 = 1.9  # U112
$ws.Cells.Item(112, 22).Value2 = 1.9  # V112
$ws.Cells.Item(112, 24).Value2 = 2  # X112
$ws.Cells.Item(112, 26).Value2 = -1  # Z112
$ws.Cells.Item(112, 27).Value2 = 0.8  # AA112
$ws.Cells.Item(112, 28).Value2 = -1  # AB112
$ws.Cells.Item(112, 29).Value2 = 0.8999999999999999  # AC112

# Row 113
$ws.Cells.Item(113, 2).Value2 = 6989521  # B113
$ws.Cells.Item(113, 6).Value2 = 'FK Radnicki Beograd'  # F113
$ws.Cells.Item(113, 7).Value2 = 'FK Indija'  # G113
$ws.Cells.Item(113, 8).Value2 = 1  # H113
$ws.Cells.Item(113, 9).Value2 = 1  # I113
$ws.Cells.Item(113, 11).Value2 = 2.4  # K113
$ws.Cells.Item(113, 13).Value2 = 2.75  # M113
$ws.Cells.Item(113, 14).Value2 = 2.875  # N113
$ws.Cells.Item(113, 15).Value2 = 2.8  # O113
$ws.Cells.Item(113, 16).Value2 = 2.45  # P113
$ws.Cells.Item(113, 17).Value2 = 0  # Q113
$ws.Cells.Item(113, 18).Value2 = 2.05  # R113
$ws.Cells.Item(113, 19).Value2 = 1.75  # S113
$ws.Cells.Item(113, 21).Value2 = 1.8  # U113
$ws.Cells.Item(113, 22).Value2 = 2  # V113
$ws.Cells.Item(113, 24).Value2 = 1.8  # X113
$ws.Cells.Item(113, 26).Value2 = 0  # Z113
$ws.Cells.Item(113, 27).Value2 = -0  # AA113
$ws.Cells.Item(113, 28).Value2 = 0  # AB113
$ws.Cells.Item(113, 29).Value2 = -0  # AC113

# Row 131
$ws.Cells.Item(131, 2).Value2 = 7497670  # B131
$ws.Cells.Item(131, 6).Value2 = 'RFK Novi Sad 1921'  # F131
$ws.Cells.Item(131, 7).Value2 = 'Radnicki Sremska Mitrovica'  # G131
$ws.Cells.Item(131, 8).Value2 = 2  # H131
$ws.Cells.Item(131, 9).Value2 = 2  # I131
$ws.Cells.Item(131, 10).Value2 = 'D'  # J131
$ws.Cells.Item(131, 11).Value2 = 3.2  # K131
$ws.Cells.Item(131, 12).Value2 = 2.9  # L131
$ws.Cells.Item(131, 13).Value2 = 2.2  # M131
$ws.Cells.Item(131, 14).Value2 = 3.3  # N131
$ws.Cells.Item(131, 15).Value2 = 2.9  # O131
$ws.Cells.Item(131, 16).Value2 = 2.15  # P131
$ws.Cells.Item(131, 18).Value2 = 1.875  # R131
$ws.Cells.Item(131, 19).Value2 = 1.925  # S131
$ws.Cells.Item(131, 21).Value2 = 1.95  # U131
$ws.Cells.Item(131, 22).Value2 = 1.85  # V131
$ws.Cells.Item(131, 24).Value2 = 1.9  # X131
$ws.Cells.Item(131, 25).Value2 = -1  # Y131
$ws.Cells.Item(131, 26).Value2 = 0.4375  # Z131
$ws.Cells.Item(131, 27).Value2 = -0.5  # AA131
$ws.Cells.Item(131, 28).Value2 = 0.95  # AB131
$ws.Cells.Item(131, 29).Value2 = -1  # AC131

# Row 133
$ws.Cells.Item(133, 2).Value2 = 7497904  # B133
$ws.Cells.Item(133, 6).Value2 = 'Sloboda Uzice'  # F133
$ws.Cells.Item(133, 7).Value2 = 'FK Indija'  # G133
$ws.Cells.Item(133, 8).Value2 = 0  # H133
$ws.Cells.Item(133, 9).Value2 = 1  # I133
$ws.Cells.Item(133, 10).Value2 = 'A'  # J133
$ws.Cells.Item(133, 11).Value2 = 2.4  # K133
$ws.Cells.Item(133, 12).Value2 = 2.875  # L133
$ws.Cells.Item(133, 13).Value2 = 2.875  # M133
$ws.Cells.Item(133, 14).Value2 = 2.875  # N133
$ws.Cells.Item(133, 15).Value2 = 3  # O133
$ws.Cells.Item(133, 16).Value2 = 2.3  # P133
$ws.Cells.Item(133, 18).Value2 = 1.75  # R133
$ws.Cells.Item(133, 19).Value2 = 2.05  # S133
$ws.Cells.Item(133, 21).Value2 = 1.925  # U133
$ws.Cells.Item(133, 22).Value2 = 1.875  # V133
$ws.Cells.Item(133, 24).Value2 = -1  # X133
$ws.Cells.Item(133, 25).Value2 = 1.3  # Y133
$ws.Cells.Item(133, 26).Value2 = -1  # Z133
$ws.Cells.Item(133, 27).Value2 = 1.05  # AA133
$ws.Cells.Item(133, 28).Value2 = -1  # AB133
$ws.Cells.Item(133, 29).Value2 = 0.875  # AC133

# Row 145
$ws.Cells.Item(145, 2).Value2 = 6989529  # B145
$ws.Cells.Item(145, 6).Value2 = 'FK Radnicki Beograd'  # F145
$ws.Cells.Item(145, 7).Value2 = 'OFK Belgrade'  # G145
$ws.Cells.Item(145, 9).Value2 = 2  # I145
$ws.Cells.Item(145, 10).Value2 = 'A'  # J145
$ws.Cells.Item(145, 11).Value2 = 4.75  # K145
$ws.Cells.Item(145, 12).Value2 = 3.75  # L145
$ws.Cells.Item(145, 13).Value2 = 1.571  # M145
$ws.Cells.Item(145, 14).Value2 = 4.75  # N145
$ws.Cells.Item(145, 15).Value2 = 3.75  # O145
$ws.Cells.Item(145, 16).Value2 = 1.571  # P145
$ws.Cells.Item(145, 17).Value2 = 1  # Q145
$ws.Cells.Item(145, 18).Value2 = 1.8  # R145
$ws.Cells.Item(145, 19).Value2 = 2  # S145
$ws.Cells.Item(145, 20).Value2 = 2.5  # T145
$ws.Cells.Item(145, 21).Value2 = 1.85  # U145
$ws.Cells.Item(145, 22).Value2 = 1.95  # V145
$ws.Cells.Item(145, 24).Value2 = -1  # X145
$ws.Cells.Item(145, 25).Value2 = 0.571  # Y145
$ws.Cells.Item(145, 26).Value2 = 0  # Z145
$ws.Cells.Item(145, 27).Value2 = -0  # AA145
$ws.Cells.Item(145, 28).Value2 = 0.8500000000000001  # AB145
$ws.Cells.Item(145, 29).Value2 = -1  # AC145

# Row 146
$ws.Cells.Item(146, 2).Value2 = 6989631  # B146
$ws.Cells.Item(146, 6).Value2 = 'Jedinstvo UB'  # F146
$ws.Cells.Item(146, 7).Value2 = 'FK Tekstilac Odzaci'  # G146
$ws.Cells.Item(146, 8).Value2 = 1  # H146
$ws.Cells.Item(146, 9).Value2 = 1  # I146
$ws.Cells.Item(146, 11).Value2 = 2.25  # K146
$ws.Cells.Item(146, 12).Value2 = 3  # L146
$ws.Cells.Item(146, 13).Value2 = 3  # M146
$ws.Cells.Item(146, 14).Value2 = 2.25  # N146
$ws.Cells.Item(146, 15).Value2 = 3  # O146
$ws.Cells.Item(146, 16).Value2 = 3  # P146
$ws.Cells.Item(146, 17).Value2 = -0.25  # Q146
$ws.Cells.Item(146, 18).Value2 = 2  # R146
$ws.Cells.Item(146, 19).Value2 = 1.8  # S146
$ws.Cells.Item(146, 20).Value2 = 2.25  # T146
$ws.Cells.Item(146, 21).Value2 = 2  # U146
$ws.Cells.Item(146, 22).Value2 = 1.8  # V146
$ws.Cells.Item(146, 24).Value2 = 2  # X146
$ws.Cells.Item(146, 26).Value2 = -0.5  # Z146
$ws.Cells.Item(146, 27).Value2 = 0.4  # AA146
$ws.Cells.Item(146, 28).Value2 = -0.5  # AB146
$ws.Cells.Item(146, 29).Value2 = 0.4  # AC146

# Row 147
$ws.Cells.Item(147, 2).Value2 = 7019002  # B147
$ws.Cells.Item(147, 6).Value2 = 'Metalac Gornji'  # F147
$ws.Cells.Item(147, 7).Value2 = 'FK Kolubara'  # G147
$ws.Cells.Item(147, 8).Value2 = 2  # H147
$ws.Cells.Item(147, 10).Value2 = 'D'  # J147
$ws.Cells.Item(147, 11).Value2 = 2.5  # K147
$ws.Cells.Item(147, 12).Value2 = 3.2  # L147
$ws.Cells.Item(147, 13).Value2 = 2.5  # M147
$ws.Cells.Item(147, 14).Value2 = 1.75  # N147
$ws.Cells.Item(147, 15).Value2 = 3.1  # O147
$ws.Cells.Item(147, 16).Value2 = 4.5  # P147
$ws.Cells.Item(147, 17).Value2 = -0.5  # Q147
$ws.Cells.Item(147, 20).Value2 = 2  # T147
$ws.Cells.Item(147, 21).Value2 = 1.9  # U147
$ws.Cells.Item(147, 22).Value2 = 1.9  # V147
$ws.Cells.Item(147, 24).Value2 = 2.1  # X147
$ws.Cells.Item(147, 25).Value2 = -1  # Y147
$ws.Cells.Item(147, 26).Value2 = -1  # Z147
$ws.Cells.Item(147, 27).Value2 = 1  # AA147
$ws.Cells.Item(147, 28).Value2 = 0.8999999999999999  # AB147

# Row 189
$ws.Cells.Item(189, 2).Value2 = 6989643  # B189
$ws.Cells.Item(189, 6).Value2 = 'OFK Vrsac'  # F189
$ws.Cells.Item(189, 7).Value2 = 'Jedinstvo UB'  # G189
$ws.Cells.Item(189, 8).Value2 = 1  # H189
$ws.Cells.Item(189, 9).Value2 = 0  # I189
$ws.Cells.Item(189, 10).Value2 = 'H'  # J189
$ws.Cells.Item(189, 12).Value2 = 3.25  # L189
$ws.Cells.Item(189, 13).Value2 = 2.75  # M189
$ws.Cells.Item(189, 14).Value2 = 2.3  # N189
$ws.Cells.Item(189, 15).Value2 = 3.3  # O189
$ws.Cells.Item(189, 16).Value2 = 2.625  # P189
$ws.Cells.Item(189, 17).Value2 = -0.25  # Q189
$ws.Cells.Item(189, 18).Value2 = 1.975  # R189
$ws.Cells.Item(189, 19).Value2 = 1.725  # S189
$ws.Cells.Item(189, 20).Value2 = 1.75  # T189
$ws.Cells.Item(189, 21).Value2 = 1.725  # U189
$ws.Cells.Item(189, 22).Value2 = 1.975  # V189
$ws.Cells.Item(189, 23).Value2 = 1.3  # W189
$ws.Cells.Item(189, 25).Value2 = -1  # Y189
$ws.Cells.Item(189, 26).Value2 = 0.9750000000000001  # Z189
$ws.Cells.Item(189, 27).Value2 = -1  # AA189
$ws.Cells.Item(189, 29).Value2 = 0.9750000000000001  # AC189

# Row 190
$ws.Cells.Item(190, 2).Value2 = 6989536  # B190
$ws.Cells.Item(190, 6).Value2 = 'FK Mladost Gat Novi Sad'  # F190
$ws.Cells.Item(190, 7).Value2 = 'FK Radnicki Beograd'  # G190
$ws.Cells.Item(190, 8).Value2 = 0  # H190
$ws.Cells.Item(190, 9).Value2 = 1  # I190
$ws.Cells.Item(190, 10).Value2 = 'A'  # J190
$ws.Cells.Item(190, 12).Value2 = 2.75  # L190
$ws.Cells.Item(190, 13).Value2 = 3.25  # M190
$ws.Cells.Item(190, 14).Value2 = 1.4  # N190
$ws.Cells.Item(190, 15).Value2 = 4  # O190
$ws.Cells.Item(190, 16).Value2 = 7  # P190
$ws.Cells.Item(190, 17).Value2 = -1.25  # Q190
$ws.Cells.Item(190, 18).Value2 = 1.875  # R190
$ws.Cells.Item(190, 19).Value2 = 1.925  # S190
$ws.Cells.Item(190, 20).Value2 = 2.5  # T190
$ws.Cells.Item(190, 21).Value2 = 1.875  # U190
$ws.Cells.Item(190, 22).Value2 = 1.925  # V190
$ws.Cells.Item(190, 23).Value2 = -1  # W190
$ws.Cells.Item(190, 25).Value2 = 6  # Y190
$ws.Cells.Item(190, 26).Value2 = -1  # Z190
$ws.Cells.Item(190, 27).Value2 = 0.925  # AA190
$ws.Cells.Item(190, 29).Value2 = 0.925  # AC190

# Row 200
$ws.Cells.Item(200, 2).Value2 = 6989337  # B200
$ws.Cells.Item(200, 6).Value2 = 'OFK Vrsac'  # F200
$ws.Cells.Item(200, 7).Value2 = 'FK Dubocica'  # G200
$ws.Cells.Item(200, 9).Value2 = 0  # I200
$ws.Cells.Item(200, 10).Value2 = 'D'  # J200
$ws.Cells.Item(200, 11).Value2 = 1.909  # K200
$ws.Cells.Item(200, 13).Value2 = 3.75  # M200
$ws.Cells.Item(200, 14).Value2 = 1.6  # N200
$ws.Cells.Item(200, 15).Value2 = 3.4  # O200
$ws.Cells.Item(200, 16).Value2 = 5.25  # P200
$ws.Cells.Item(200, 17).Value2 = -0.75  # Q200
$ws.Cells.Item(200, 18).Value2 = 1.95  # R200
$ws.Cells.Item(200, 24).Value2 = 2.4  # X200
$ws.Cells.Item(200, 25).Value2 = -1  # Y200

# Row 201
$ws.Cells.Item(201, 2).Value2 = 6989717  # B201
$ws.Cells.Item(201, 6).Value2 = 'FK Indija'  # F201
$ws.Cells.Item(201, 7).Value2 = 'Metalac Gornji'  # G201
$ws.Cells.Item(201, 9).Value2 = 1  # I201
$ws.Cells.Item(201, 10).Value2 = 'A'  # J201
$ws.Cells.Item(201, 11).Value2 = 1.833  # K201
$ws.Cells.Item(201, 13).Value2 = 4  # M201
$ws.Cells.Item(201, 14).Value2 = 2  # N201
$ws.Cells.Item(201, 15).Value2 = 3  # O201
$ws.Cells.Item(201, 16).Value2 = 3.5  # P201
$ws.Cells.Item(201, 17).Value2 = -0.5  # Q201
$ws.Cells.Item(201, 18).Value2 = 2.05  # R201
$ws.Cells.Item(201, 24).Value2 = -1  # X201
$ws.Cells.Item(201, 25).Value2 = 2.5  # Y201

